$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $cellRef, $val) {
    $c = $range.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "27.176.98"
$ws.Range("E2").Value = "  +0.47%  "
Set-TextValue $ws "D3" "1.685.31"
$ws.Range("E3").Value = "  +0.17%  "
Set-TextValue $ws "D5" "215.94"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("E7").Value = "  +0.16%  "
Set-TextValue $ws "D8" "23.13"
$ws.Range("E8").Value = "  +7.86%  "
Set-TextValue $ws "D9" "0.260"
$ws.Range("E9").Value = "  +3.38%  "
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("E11").Value = "  +0.25%  "
Set-TextValue $ws "D12" "1.923.96"
$ws.Range("E12").Value = "  +0.21%  "
Set-TextValue $ws "D13" "1.693.78"
$ws.Range("E13").Value = "  +0.77%  "
Set-TextValue $ws "D14" "4.19"
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("E15").Value = "  +3.90%  "
Set-TextValue $ws "D16" "66.94"
$ws.Range("E16").Value = "  +1.24%  "
Set-TextValue $ws "D17" "27.182.07"
$ws.Range("E17").Value = "  +0.48%  "
Set-TextValue $ws "D18" "236.07"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("E19").Value = "  -2.16%  "
Set-TextValue $ws "D20" "0.0₃0743"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("E21").Value = "  +0.21%  "
Set-TextValue $ws "D22" "4.56"
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("E23").Value = "  +3.93%  "
$ws.Range("E24").Value = "  -2.86%  "
Set-TextValue $ws "D25" "147.19"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  +1.13%  "
Set-TextValue $ws "D27" "16.43"
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("E29").Value = "  +0.22%  "
Set-TextValue $ws "D30" "0.0505"
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  +1.26%  "
Set-TextValue $ws "D33" "1.544.90"
$ws.Range("E33").Value = "  +1.81%  "
Set-TextValue $ws "D34" "3.24"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("E35").Value = "  -1.40%  "
Set-TextValue $ws "D36" "0.604"
$ws.Range("E36").Value = "  +2.41%  "
Set-TextValue $ws "D37" "0.947"
$ws.Range("E37").Value = "  +3.07%  "
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("E44").Value = "  -0.77%  "
Set-TextValue $ws "D45" "1.831.79"
$ws.Range("E45").Value = "  +0.35%  "
Set-TextValue $ws "D46" "0.789"
$ws.Range("E46").Value = "  +0.91%  "
Set-TextValue $ws "D47" "90.25"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E48").Value = "  +5.32%  "
Set-TextValue $ws "D49" "1.61"
Set-TextValue $ws "D50" "8.32"
$ws.Range("E50").Value = "  +5.96%  "
$ws.Range("E51").Value = "  -0.76%  "
